$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark that currently sits
#     right after the "IN FUTURO" run (before ": generalizzare ...").
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# --- Step 2: merge "IN FUTURO" and ": generalizzare " (which used to be
#     two separate runs split by the bookmark) into a single run of text
#     "IN FUTURO: generalizzare " by doing a Find/Replace over the
#     now-contiguous text.
$d.Content.Find.Execute("IN FUTURO: generalizzare ", $true, $false, $false, $false, $false, $true, 1, $false, "IN FUTURO: generalizzare ", 2)

# --- Step 3: append " (FATTO)" right after "...a tutti gli esempi interattivi"
#     in the "MARCO: aggiungere il tag interactive..." paragraph, then move
#     the "_GoBack" bookmark there (after the newly appended text).
$r = $d.Content
$r.Find.Execute("a tutti gli esempi interattivi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

# Insert the new text plus a sacrificial trailing character. Creating a
# zero-length bookmark exactly at the last character position of a
# paragraph (right before its paragraph mark) is unreliable, so we give
# ourselves a safe interior position to drop the bookmark at, then trim
# the sacrificial character off afterwards.
$r.InsertAfter(" (FATTO)X")
$r.LanguageID = "it-IT"

$bmPos = $r.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sacrifice = $d.Range($r.End - 1, $r.End)
$sacrifice.Delete()
